$d = $word.ActiveDocument

# Rename the rule tag: @BEGIN BodyTextFontIncorrect -> @BEGIN TextFontIncorrect
$d.Content.Find.Execute("@BEGIN BodyTextFontIncorrect", $true, $false, $false, $false, $false,
                         $true, 1, $false, "@BEGIN TextFontIncorrect", 2)

# Update the rule's body text: drop "работы" so it applies to all text
$d.Content.Find.Execute("Текст работы должен быть написан в фонте ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Текст должен быть написан в фонте ", 2)
